# fix: add domain field in EU funds report (#13)
#
# Adds a new "Domain" column (M) at the end of the header row on the
# "Data" sheet, extends the AutoFilter / _FilterDatabase range from
# A1:L1 to A1:M1, and updates column widths so that the new last
# column (M) takes over the wider formatting that used to belong to
# the old last column (L), while K:L share the narrower width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Copy the formatting (style) of the current last header cell (L1,
# "Postal Address") onto the new header cell M1, then set its text.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Domain"

# Clear the clipboard marquee / copy mode.
$excel.CutCopyMode = 0

# Re-create the AutoFilter over the expanded A1:M1 range.
$ws.AutoFilterMode = $false
$ws.Range("A1:M1").AutoFilter() | Out-Null

# Column widths: M (new last column) takes the width that L used to
# have (22.5); K and L now share the narrower width (18.5).
$oldK = $ws.Columns.Item(11).ColumnWidth
$oldL = $ws.Columns.Item(12).ColumnWidth
$ws.Columns.Item(13).ColumnWidth = $oldL
$ws.Columns.Item(12).ColumnWidth = $oldK

# Update the hidden _FilterDatabase defined name to match the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$M`$1"
    }
}

# Match the author's final cursor position / selection.
[void]$ws.Range("J10").Select()
